$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B:E values (regenerated s_vals), G recomputed as sum of B:E for each row.
$data = @{
    2  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    3  = @(0.7287194209349384, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569)
    4  = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    5  = @(0.1554434735375247, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
    6  = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    7  = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 6.48142807727062)
    8  = @(0.3464964993005633, 0.05231270169004087, 0.1529057820181812, 0.4998867070740569)
    9  = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569)
    10 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    11 = @(1.505614041169197, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569)
    12 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
}

$gvals = @{
    2  = 4.371470058157054
    3  = 2.27892381503245
    4  = 3.594575437922795
    5  = 5.39116606621397
    6  = 3.594575437922795
    7  = 10.35301142835362
    8  = 1.051601690082842
    9  = 5.488907176552729
    10 = 6.048734245549538
    11 = 2.495991366269901
    12 = 4.371470058157054
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $gvals[$row]
}
